$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 98) following the same layout as the
# existing rows: column A holds the date as plain text, columns B-E
# hold numeric values.
$ws.Cells.Item(98, 1).Value = "2024-11-08 00:00:00"
$ws.Cells.Item(98, 2).Value = 75400
$ws.Cells.Item(98, 3).Value = 10530.73
$ws.Cells.Item(98, 4).Value = 9319.23
$ws.Cells.Item(98, 5).Value = 7.1592
